# Update "想去人数" (want-to-go count) figures on the "展览" and "全部类型" sheets
# to reflect the latest scrape (gh-pages output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 3464   # was 3462
$wsExpo.Range("F3").Value = 32     # was 31
$wsExpo.Range("F5").Value = 1921   # was 1912
$wsExpo.Range("F6").Value = 138    # was 137
$wsExpo.Range("F7").Value = 347    # was 345

# --- Sheet "全部类型" (All types) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 3464    # was 3462
$wsAll.Range("F3").Value = 32      # was 31
$wsAll.Range("F5").Value = 1921    # was 1912
$wsAll.Range("F6").Value = 138     # was 137
$wsAll.Range("F8").Value = 347     # was 345
